# Update the "2025" year-to-date standings rows (201-210) on Sheet1 with
# refreshed figures (the source data was re-pulled from Feb to May).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 201 - Andy stays rank 1, figures refreshed
$ws.Range("D201").Value = 26
$ws.Range("F201").Value = 26
$ws.Range("G201").Value = 94150
$ws.Range("H201").Value = 110
$ws.Range("I201").Value = 60

# Row 202 - Prashant stays rank 2, figures refreshed
$ws.Range("D202").Value = 22
$ws.Range("F202").Value = 22
$ws.Range("G202").Value = 67150
$ws.Range("H202").Value = 80
$ws.Range("I202").Value = 30

# Row 203 - now Pepe (rank 3)
$ws.Range("B203").Value = "Pepe"
$ws.Range("D203").Value = 21
$ws.Range("F203").Value = 21
$ws.Range("G203").Value = 67550
$ws.Range("H203").Value = 30
$ws.Range("I203").Value = -20
$ws.Range("K203").Value = 364

# Row 204 - now Matt (rank 4)
$ws.Range("B204").Value = "Matt"
$ws.Range("D204").Value = 21
$ws.Range("F204").Value = 21
$ws.Range("G204").Value = 64700
$ws.Range("H204").Value = 60
$ws.Range("I204").Value = 10
$ws.Range("K204").Value = 362

# Row 205 - now Richard (rank 5)
$ws.Range("B205").Value = "Richard"
$ws.Range("D205").Value = 19
$ws.Range("F205").Value = 19
$ws.Range("G205").Value = 67850
$ws.Range("H205").Value = 100
$ws.Range("I205").Value = 50
$ws.Range("K205").Value = 366

# Row 206 - Maisy stays rank 6, figures refreshed
$ws.Range("D206").Value = 14
$ws.Range("F206").Value = 14
$ws.Range("G206").Value = 47400
$ws.Range("H206").Value = 20

# Row 207 - now Jon (rank 7)
$ws.Range("B207").Value = "Jon"
$ws.Range("D207").Value = 7
$ws.Range("F207").Value = 7
$ws.Range("G207").Value = 30950
$ws.Range("H207").Value = 0
$ws.Range("I207").Value = -50
$ws.Range("K207").Value = 357

# Row 208 - now Mark (rank 8)
$ws.Range("B208").Value = "Mark"
$ws.Range("D208").Value = 7
$ws.Range("F208").Value = 7
$ws.Range("G208").Value = 26450
$ws.Range("H208").Value = 10
$ws.Range("I208").Value = -30
$ws.Range("K208").Value = 361

# Row 209 - Anthony stays rank 9, figures refreshed
$ws.Range("D209").Value = 7
$ws.Range("F209").Value = 7
$ws.Range("G209").Value = 21000
$ws.Range("I209").Value = -10

# Row 210 - Alex stays rank 10, figures refreshed
$ws.Range("D210").Value = 6
$ws.Range("F210").Value = 6
$ws.Range("G210").Value = 30950
$ws.Range("H210").Value = 10
$ws.Range("I210").Value = -20
